$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rows 23 downward are being re-shuffled/expanded to make room for four new
# "decimal ...InUsd"/devaluation member rows. Old content occupied rows
# 23:46 (one value per row, in column C or D). The new layout occupies
# rows 23:50.
#
# New shared-string table entries get appended in the order their text is
# first assigned to a cell, so write the four brand-new strings first (in
# the same order they appear in the shared string table of the target
# workbook), and only then fill in the rest of the block, re-using the
# exact original text for the untouched rows so Excel resolves them back
# to the very same shared-string entries (including the two rich-text
# "List<...>" entries).
$ws.Range("C41").Value = "          decimal TotalPercentInUsd"
$ws.Range("C49").Value = "          decimal EstimatedDevaluationInUsd"
$ws.Range("C48").Value = "          decimal EstimatedCurrencyRateOnFinish"
$ws.Range("C45").Value = "          decimal CurrentDevaluationInUsd"

$rows = @(
    @(23, "C", "          DepositStates State"),
    @(24, "C", "          List<DepositDailyLine> DailyTable"),
    @(25, "D", "DateTime Date"),
    @(26, "D", "decimal Balance"),
    @(27, "D", "decimal DepoRate"),
    @(28, "D", "decimal DayProcents"),
    @(29, "D", "decimal NotPaidProcents"),
    @(30, "D", "decimal CurrencyRate"),
    @(31, "D", "decimal DayDevaluation"),
    @(32, "C", "          List<DepositTransaction> Traffic"),
    @(33, "D", "DateTime Timestamp"),
    @(34, "D", "TransactionType"),
    @(35, "D", "decimal Amount"),
    @(36, "D", "Currency"),
    @(37, "D", "decimal AmountInUsd"),
    @(38, "D", "string Comment"),
    @(39, "C", "          decimal TotalMyIns"),
    @(40, "C", "          decimal TotalPercent"),
    @(42, "C", "          decimal TotalMyOuts"),
    @(43, "C", "          decimal CurrentBalance"),
    @(44, "C", "          decimal CurrentProfitInUsd"),
    @(46, "C", "          decimal EstimatedProcentsInThisMonth"),
    @(47, "C", "          decimal EstimatedProcents"),
    @(50, "C", "          decimal EstimatedProfitInUsd")
)

foreach ($row in $rows) {
    $r = $row[0]
    $col = $row[1]
    $text = $row[2]
    $ws.Range("$col$r").Value = $text
}

# Clear the stale cells left over from the old layout whose column
# assignment changed (old row 24 was D, old row 31 was C, old row 32 was D,
# old row 38 was C).
$ws.Range("D24").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("C38").ClearContents()

$ws.Range("F9").Select()
